$d = $word.ActiveDocument

# --- Step 1: split the trailing " task" run of the "My 4th task" paragraph
# into a plain " " run followed by a "task" run wrapped in gramStart/gramEnd
# proofing-error markers (mirrors the existing "own" paragraph above it).
$p4 = $d.Paragraphs(4)
$rng = $d.Range($p4.Range.End - 6, $p4.Range.End - 1)
$xmlSplit = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>task</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xmlSplit)

# --- Step 2: append a new "My 5th task" paragraph (same run/formatting shape
# as the 4th one, but without proofing markers) right after it, followed by a
# new blank paragraph.
$p4 = $d.Paragraphs(4)
$insertPoint = $d.Range($p4.Range.End, $p4.Range.End)
$xmlNew = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">My </w:t></w:r><w:r><w:t>5</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>th</w:t></w:r><w:r><w:t xml:space="preserve"> task</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($xmlNew)
